# Apply market-price / leve-profit data refresh pulled from the scheduled Sheets runner.
# For every touched row: columns H (currentAveragePrice), I (currentAveragePriceNQ),
# J (currentAveragePriceHQ), K (LevePriceNQ), L (LevePriceHQ), M (LeveProfitNQ) and
# N (LeveProfitHQ) are refreshed with the latest computed figures. A few rows gain or
# lose a cell entirely when a profit figure becomes newly computable / no longer applies.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 28 (Leve Item ID 27772)
$ws.Range("H28").Value = 1188
$ws.Range("I28").Value = 202.18182
$ws.Range("K28").Value = 202.18182
$ws.Range("M28").Value = 282.81818
# Row 41 (Leve Item ID 5478)
$ws.Range("H41").Value = 1046.5
$ws.Range("I41").Value = 1089.2667
$ws.Range("J41").Value = 832.6667
$ws.Range("K41").Value = 1089.2667
$ws.Range("L41").Value = 832.6667
$ws.Range("M41").Value = -649.2666999999999
$ws.Range("N41").Value = -1712.6667
# Row 86 (Leve Item ID 12603)
$ws.Range("H86").Value = 8252.857
$ws.Range("I86").Value = 7943.75
$ws.Range("J86").Value = 8665
$ws.Range("K86").Value = 7943.75
$ws.Range("L86").Value = 8665
$ws.Range("M86").Value = -6820.75
$ws.Range("N86").Value = -10911
# Row 89 (Leve Item ID 12603)
$ws.Range("H89").Value = 8252.857
$ws.Range("I89").Value = 7943.75
$ws.Range("J89").Value = 8665
$ws.Range("K89").Value = 39718.75
$ws.Range("L89").Value = 43325
$ws.Range("M89").Value = -34102.75
$ws.Range("N89").Value = -54557
# Row 92 (Leve Item ID 19901)
$ws.Range("H92").Value = 1122
$ws.Range("I92").Value = 1294.2354
$ws.Range("K92").Value = 1294.2354
$ws.Range("M92").Value = -46.23540000000003
# Row 107 (Leve Item ID 27766)
$ws.Range("H107").Value = 580.6087
$ws.Range("I107").Value = 697.94116
$ws.Range("J107").Value = 248.16667
$ws.Range("K107").Value = 697.94116
$ws.Range("L107").Value = 248.16667
$ws.Range("M107").Value = 1222.05884
$ws.Range("N107").Value = -4088.16667
# Row 112 (Leve Item ID 27960)
$ws.Range("H112").Value = 64578.117
$ws.Range("J112").Value = 64578.117
$ws.Range("L112").Value = 193734.351
$ws.Range("N112").Value = -195950.351

$ws = $wb.Worksheets.Item("ARM")
# Row 2 (Leve Item ID 27713)
$ws.Range("H2").Value = 830.25
$ws.Range("J2").Value = 918.375
$ws.Range("L2").Value = 918.375
$ws.Range("N2").Value = -1144.375
# Row 32 (Leve Item ID 44147)
$ws.Range("H32").Value = 7847.1
$ws.Range("I32").Value = 3228.7327
$ws.Range("K32").Value = 3228.7327
$ws.Range("M32").Value = -2941.7327
# Row 45 (Leve Item ID 27714)
$ws.Range("H45").Value = 8690.333000000001
$ws.Range("I45").Value = 10941.363
$ws.Range("K45").Value = 10941.363
$ws.Range("M45").Value = -10564.363
# Row 97 (Leve Item ID 19941)
$ws.Range("H97").Value = 2430.7
$ws.Range("I97").Value = 2111
$ws.Range("J97").Value = 3389.8
$ws.Range("K97").Value = 2111
$ws.Range("L97").Value = 3389.8
$ws.Range("M97").Value = -1615
$ws.Range("N97").Value = -4381.8
# Row 110 (Leve Item ID 27708)
$ws.Range("H110").Value = 4789.3716
$ws.Range("J110").Value = 3163.7058
$ws.Range("L110").Value = 3163.7058
$ws.Range("N110").Value = -7253.7058
# Row 116 (Leve Item ID 27713)
$ws.Range("H116").Value = 830.25
$ws.Range("J116").Value = 918.375
$ws.Range("L116").Value = 918.375
$ws.Range("N116").Value = -5506.375

$ws = $wb.Worksheets.Item("BSM")
# Row 3 (Leve Item ID 27713)
$ws.Range("H3").Value = 830.25
$ws.Range("J3").Value = 918.375
$ws.Range("L3").Value = 918.375
$ws.Range("N3").Value = -1146.375
# Row 22 (Leve Item ID 5092)
$ws.Range("H22").Value = 348264.4
$ws.Range("J22").Value = 483424.8
$ws.Range("L22").Value = 483424.8
$ws.Range("N22").Value = -483770.8
# Row 94 (Leve Item ID 19939)
$ws.Range("H94").Value = 1130.9231
$ws.Range("I94").Value = 1192.6522
$ws.Range("K94").Value = 1192.6522
$ws.Range("M94").Value = -741.6522
# Row 99 (Leve Item ID 19943)
$ws.Range("H99").Value = 2849.0981
$ws.Range("I99").Value = 2101.8333
$ws.Range("J99").Value = 3916.6191
$ws.Range("K99").Value = 2101.8333
$ws.Range("L99").Value = 3916.6191
$ws.Range("M99").Value = -603.8332999999998
$ws.Range("N99").Value = -6912.6191
# Row 105 (Leve Item ID 19947)
$ws.Range("H105").Value = 5827.1333
$ws.Range("I105").Value = 5986.0835
$ws.Range("K105").Value = 5986.0835
$ws.Range("M105").Value = -4239.0835
# Row 118 (Leve Item ID 27137)
$ws.Range("H118").Value = 35555
$ws.Range("J118").Value = 35555
$ws.Range("L118").Value = 35555
$ws.Range("N118").Value = -38869
# Row 134 (Leve Item ID 43998)
$ws.Range("H134").Value = 1424.8889
$ws.Range("I134").Value = 1382.88
$ws.Range("J134").Value = 1950
$ws.Range("K134").Value = 4148.64
$ws.Range("L134").Value = 5850
$ws.Range("M134").Value = -1613.64
$ws.Range("N134").Value = -10920

$ws = $wb.Worksheets.Item("CRP")
# Row 31 (Leve Item ID 44023)
$ws.Range("H31").Value = 2667.182
$ws.Range("J31").Value = 3349.25
$ws.Range("L31").Value = 3349.25
$ws.Range("N31").Value = -3939.25
# Row 34 (Leve Item ID 44023)
$ws.Range("H34").Value = 2667.182
$ws.Range("J34").Value = 3349.25
$ws.Range("L34").Value = 3349.25
$ws.Range("N34").Value = -3753.25
# Row 62 (Leve Item ID 12580)
$ws.Range("H62").Value = 2307.7144
$ws.Range("I62").Value = 2400.8
$ws.Range("J62").Value = 2075
$ws.Range("K62").Value = 2400.8
$ws.Range("L62").Value = 2075
$ws.Range("M62").Value = -1776.8
$ws.Range("N62").Value = -3323
# Row 65 (Leve Item ID 12580)
$ws.Range("H65").Value = 2307.7144
$ws.Range("I65").Value = 2400.8
$ws.Range("J65").Value = 2075
$ws.Range("K65").Value = 12004
$ws.Range("L65").Value = 10375
$ws.Range("M65").Value = -8884
$ws.Range("N65").Value = -16615

$ws = $wb.Worksheets.Item("CUL")
# Row 4 (Leve Item ID 4650)
$ws.Range("H4").Value = 2316496
$ws.Range("I4").Value = 2279645.5
$ws.Range("J4").Value = 2500748.5
$ws.Range("K4").Value = 6838936.5
$ws.Range("L4").Value = 7502245.5
$ws.Range("M4").Value = -6838824.5
$ws.Range("N4").Value = -7502469.5
# Row 62 (Leve Item ID 12845)
$ws.Range("H62").Value = 4940.1
$ws.Range("I62").Value = 950.5
$ws.Range("K62").Value = 2851.5
$ws.Range("M62").Value = -2165.5
# Row 65 (Leve Item ID 12845)
$ws.Range("H65").Value = 4940.1
$ws.Range("I65").Value = 950.5
$ws.Range("K65").Value = 8554.5
$ws.Range("M65").Value = -5122.5
# Row 132 (Leve Item ID 43972)
$ws.Range("H132").Value = 1776.3077
$ws.Range("I132").Value = 1413.0555
$ws.Range("J132").Value = 2593.625
$ws.Range("K132").Value = 12717.4995
$ws.Range("L132").Value = 23342.625
$ws.Range("M132").Value = -10187.4995
$ws.Range("N132").Value = -28402.625
# Row 139 (Leve Item ID 44102)
$ws.Range("H139").Value = 2866.3333
$ws.Range("I139").Value = 2749.5
$ws.Range("J139").Value = 3100
$ws.Range("K139").Value = 8248.5
$ws.Range("L139").Value = 9300
$ws.Range("M139").Value = -3108.5
$ws.Range("N139").Value = -19580

$ws = $wb.Worksheets.Item("GSM")
# Row 97 (Leve Item ID 19940)
$ws.Range("H97").Value = 1897
$ws.Range("I97").Value = 2163
$ws.Range("K97").Value = 2163
$ws.Range("M97").Value = -1667
# Row 98 (Leve Item ID 18359)
$ws.Range("H98").Value = 0
$ws.Range("J98").Value = 0
$ws.Range("L98").Value = 0
$ws.Range("N98").ClearContents()
# Row 105 (Leve Item ID 18671)
$ws.Range("H105").Value = 40000
$ws.Range("J105").Value = 40000
$ws.Range("L105").Value = 40000
$ws.Range("N105").Value = -46988
# Row 132 (Leve Item ID 44008)
$ws.Range("H132").Value = 2460.9
$ws.Range("I132").Value = 1634.1818
$ws.Range("J132").Value = 4734.375
$ws.Range("K132").Value = 4902.5454
$ws.Range("L132").Value = 14203.125
$ws.Range("M132").Value = -2372.5454
$ws.Range("N132").Value = -19263.125
# Row 137 (Leve Item ID 43226)
$ws.Range("H137").Value = 100780
$ws.Range("J137").Value = 100780
$ws.Range("L137").Value = 100780
$ws.Range("N137").Value = -110980

$ws = $wb.Worksheets.Item("LTW")
# Row 22 (Leve Item ID 5277)
$ws.Range("H22").Value = 1854.8918
$ws.Range("I22").Value = 2160.3572
$ws.Range("K22").Value = 2160.3572
$ws.Range("M22").Value = -1865.3572
# Row 27 (Leve Item ID 5277)
$ws.Range("H27").Value = 1854.8918
$ws.Range("I27").Value = 2160.3572
$ws.Range("K27").Value = 2160.3572
$ws.Range("M27").Value = -2053.3572
# Row 93 (Leve Item ID 19993)
$ws.Range("H93").Value = 3037.0908
$ws.Range("I93").Value = 3200.6667
$ws.Range("J93").Value = 2840.8
$ws.Range("K93").Value = 3200.6667
$ws.Range("L93").Value = 2840.8
$ws.Range("M93").Value = -1952.6667
$ws.Range("N93").Value = -5336.8
# Row 100 (Leve Item ID 19995)
$ws.Range("H100").Value = 402897
$ws.Range("J100").Value = 503121.25
$ws.Range("L100").Value = 503121.25
$ws.Range("N100").Value = -504203.25
# Row 123 (Leve Item ID 35408)
$ws.Range("H123").Value = 35277
$ws.Range("J123").Value = 35277
$ws.Range("L123").Value = 35277
$ws.Range("N123").Value = -45077
# Row 134 (Leve Item ID 42024)
$ws.Range("H134").Value = 125000
$ws.Range("J134").Value = 125000
$ws.Range("L134").Value = 125000
$ws.Range("N134").Value = -135140
# Row 137 (Leve Item ID 43296)
$ws.Range("H137").Value = 103553.625
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 103553.625
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 103553.625
$ws.Range("M137").ClearContents()
$ws.Range("N137").Value = -113753.625

$ws = $wb.Worksheets.Item("WVR")
# Row 4 (Leve Item ID 2996)
$ws.Range("H4").Value = 3597.3572
$ws.Range("I4").Value = 3000
$ws.Range("J4").Value = 3643.3076
$ws.Range("K4").Value = 3000
$ws.Range("L4").Value = 3643.3076
$ws.Range("M4").Value = -2887
$ws.Range("N4").Value = -3869.3076
# Row 5 (Leve Item ID 3515)
$ws.Range("H5").Value = 4500
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 4500
$ws.Range("K5").Value = 0
$ws.Range("L5").Value = 4500
$ws.Range("M5").ClearContents()
$ws.Range("N5").Value = -4724
# Row 127 (Leve Item ID 35414)
$ws.Range("H127").Value = 41510
$ws.Range("J127").Value = 41510
$ws.Range("L127").Value = 41510
$ws.Range("N127").Value = -51430
# Row 136 (Leve Item ID 44031)
$ws.Range("H136").Value = 1530.4762
$ws.Range("J136").Value = 3750
$ws.Range("L136").Value = 11250
$ws.Range("N136").Value = -16350
